# Add a new "Uridine" entry to the inclusion list and re-sort the table
# by RT_PH (column B), matching the author's re-upload of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right before the current "Inosine" row (row 11) and
# populate it with the new analyte's data.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = "Uridine"
$ws.Cells.Item(11, 2).Value = 243.0623
$ws.Cells.Item(11, 3).Value = 31

# Re-sort the data (A2:C13) ascending by RT_PH (column B), same as Excel's
# Data > Sort feature, so the new row lands in its correct position.
$dataRange = $ws.Range("A2:C13")
$sortKey = $ws.Range("B2:B13")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Reflect the selection left active after performing the sort.
$ws.Range("A2:C13").Select()
